# Apply odds updates for the FlashScore 2024-10-22 weekly games workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (KfdOaoV7)
$ws.Range("BD3").Value = 151

# Row 4 (CjY4N3UI)
$ws.Range("O4").Value = 1.17
$ws.Range("P4").Value = 5
$ws.Range("Q4").Value = 1.57
$ws.Range("R4").Value = 2.38

# Row 8 (6Dgvx3qC)
$ws.Range("N8").Value = 10
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.8

# Row 9 (Ao0VwPEa)
$ws.Range("G9").Value = 2.9
$ws.Range("I9").Value = 2.4
$ws.Range("J9").Value = 3.5
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 10
$ws.Range("AA9").Value = 23
$ws.Range("AH9").Value = 8.5
$ws.Range("AM9").Value = 26
$ws.Range("AN9").Value = 5
$ws.Range("AY9").Value = 21

# Row 32 (OUJSXR8e)
$ws.Range("N32").Value = 17

# Row 39 (8ni80y74)
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 3.75
$ws.Range("I39").Value = 1.75
$ws.Range("J39").Value = 4.5
$ws.Range("L39").Value = 2.38
$ws.Range("AA39").Value = 34
$ws.Range("AH39").Value = 8
$ws.Range("AI39").Value = 9
$ws.Range("AL39").Value = 13
$ws.Range("AP39").Value = 29
$ws.Range("AX39").Value = 9
$ws.Range("BA39").Value = 41
